$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.46580998
$ws.Range("H2").Value = 169.88336374
$ws.Range("M2").Value = 0.2575856796067603
$ws.Range("N2").Value = 35.65025453032664
$ws.Range("G3").Value = 1.51946332
$ws.Range("H3").Value = 171.45638814
$ws.Range("M3").Value = 0.4563816235617154
$ws.Range("N3").Value = 60.62167444470195
$ws.Range("G4").Value = 1.31476032
$ws.Range("H4").Value = 60.33712061999999
$ws.Range("M4").Value = 0.2901926646659931
$ws.Range("N4").Value = 22.75476803659272
$ws.Range("G5").Value = 0.7536861799999999
$ws.Range("H5").Value = 48.68516868
$ws.Range("M5").Value = 0.2075764034850184
$ws.Range("N5").Value = 17.96604165128577
$ws.Range("G6").Value = 0.5826141199999999
$ws.Range("H6").Value = 16.24644572
$ws.Range("M6").Value = 0.1705166284447882
$ws.Range("N6").Value = 7.710043693630528
$ws.Range("G7").Value = 0.3537025600000001
$ws.Range("H7").Value = 12.07943624
$ws.Range("M7").Value = 0.1064298320441421
$ws.Range("N7").Value = 4.768850306221507
$ws.Range("G8").Value = 0.42461368
$ws.Range("H8").Value = 4.05563942
$ws.Range("M8").Value = 0.1572723374045257
$ws.Range("N8").Value = 2.259023269990147
$ws.Range("G9").Value = 0.47369198
$ws.Range("H9").Value = 8.272052520000001
$ws.Range("M9").Value = 0.1169436231843642
$ws.Range("N9").Value = 3.776438336769833
$ws.Range("G10").Value = 0.2038068
$ws.Range("H10").Value = 1.46559084
$ws.Range("M10").Value = 0.08619550150860913
$ws.Range("N10").Value = 0.9145212355313407
$ws.Range("G11").Value = 0.26120484
$ws.Range("H11").Value = 3.65293844
$ws.Range("M11").Value = 0.08015625125843905
$ws.Range("N11").Value = 2.300968590328717
$ws.Range("G12").Value = 0.13002958
$ws.Range("H12").Value = 0.80684384
$ws.Range("M12").Value = 0.06834593545933235
$ws.Range("N12").Value = 0.6244658429654223
$ws.Range("G13").Value = 0.14244442
$ws.Range("H13").Value = 1.61395626
$ws.Range("M13").Value = 0.04614002941581838
$ws.Range("N13").Value = 0.940494734128831
